$wb = $excel.ActiveWorkbook
$wsTestCases = $wb.Worksheets.Item("TestCases")
$wsTestData = $wb.Worksheets.Item("TestData")

# Update the "Runmode" flag for the OpenAccountTest row on TestData from Y to N
$wsTestData.Range("A4").Value = "N"
$wsTestData.Activate()
$wsTestData.Range("A4").Select() | Out-Null

# Update the "Runmode" flag for the OpenAccountTest row on TestCases from Y to N
$wsTestCases.Range("B3").Value = "N"

# Make TestCases the active/selected sheet, with B3 selected
$wsTestCases.Activate()
$wsTestCases.Range("B3").Select() | Out-Null
